$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing hour values (rows 6, 13, 15) ---
$ws.Range("C6").Value = 1.5
$ws.Range("C13").Value = 5
$ws.Range("C15").Value = 2.5

# --- Fill in the new journal rows (17-29) ---

# Row 17
$ws.Range("A17").Value = 43178
$ws.Range("B17").Value = 'recherches sur le conrôleur de calque (mise en place d''un hiérarche de classes)'
$ws.Range("C17").NumberFormat = "0.0"
$ws.Range("C17").Value = 1.5

# Row 18
$ws.Range("A18").Value = 43192
$ws.Range("B18").Value = "discussions sur l'implémentation de la sauvegarde + mise au point sur la suite du projet"
$ws.Range("C18").Value = 1.5

# Row 19
$ws.Range("A19").Value = 43199
$ws.Range("B19").Value = 'recherche sur la modélisation du contrôler'
$ws.Range("C19").Value = 1.5

# Row 20
$ws.Range("A20").Value = 43205
$ws.Range("B20").Value = 'refactorisation du code du programme, implémentation du pencil'
$ws.Range("C20").Value = 9

# Row 21 (Heures column stored as the text "0.25")
$ws.Range("A21").Value = 43206
$ws.Range("B21").Value = 'Présentation intérmédiaire du projet'
$ws.Range("C21").Value = 0.25

# Row 22
$ws.Range("A22").Value = 43219
$ws.Range("B22").Value = "Restructuration du projet, mise en place correcte du pencil et debut d'implémentation de la gomme"
$ws.Range("C22").Value = 8

# Row 23
$ws.Range("A23").Value = 43220
$ws.Range("B23").Value = "Recherches sur l'implémentation de la gomme (beaucou de difficultés à faire une gomme tirant un trait continu)"
$ws.Range("C23").NumberFormat = "0.00"
$ws.Range("C23").Value = 1.5

# Row 24
$ws.Range("A24").Value = 43230
$ws.Range("B24").Value = 'Implémentation des formes'
$ws.Range("C24").Value = 13

# Row 25
$ws.Range("A25").Value = 43231
$ws.Range("B25").Value = 'Implémentation des formes + rédaction du rapport'
$ws.Range("C25").Value = 10

# Row 26
$ws.Range("A26").Value = 43234
$ws.Range("B26").Value = 'Rédaction du rapport + recherche sur un moyen de rendre la base des traits de crayon ronde (jusque là, le trait de crayon dessinait des suites de carrés)'
$ws.Range("C26").NumberFormat = "0.00"
$ws.Range("C26").Value = 1.5

# Row 27
$ws.Range("A27").Value = 43237
$ws.Range("B27").Value = 'Rédaction du rapport'
$ws.Range("C27").Value = 1.5

# Row 28
$ws.Range("A28").Value = 43239
$ws.Range("B28").Value = "Rédaction du rapport et du manuel d'utilisation"
$ws.Range("C28").Value = 5

# Row 29
$ws.Range("A29").Value = 43240
$ws.Range("B29").Value = 'Redaction du rapport et du manuel d''utilisation '
$ws.Range("C29").Value = 8

# --- Set row heights to match the wrapped-text layout ---
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(22).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 45

# --- Remove the now-unused trailing blank rows (30-32); this also shifts the
#     Total row up from 33 to 30 and rewrites its SUM formula range. ---
$ws.Range("A30:C32").Delete()

# --- Update the active selection to match the author's last editing position ---
$ws.Range("B35").Select()
